$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "Listeria_monocytogenes_1408MLGX"

# Keep the (broken) defined name's sheet-qualified reference in sync with the new sheet name
$wb.Names.Item(1).RefersTo = "=Listeria_monocytogenes_1408MLGX!#REF!"

# Update checksum values (sha256sumAssembly column, H) to reflect nonrefseq loci
$ws.Range("H9").Value = "d531614e1b1ba70261a6840d0ca422df90a963d6b9b696d635b8bbcc05cc07af"
$ws.Range("H12").Value = "aa128be8a56a30c5d4a55196b341c34fb73664e741723417d573f6f338691f15"
$ws.Range("H14").Value = "f4c5606eadc763af8e1b284f7eb81ac3414eedc36b5abe4bcf54261e356bc6d4"
$ws.Range("H20").Value = "5727b8d883ef850863598d63f533f3e53c68be0e8fa82c06c9e81e768ba251a4"
$ws.Range("H25").Value = "f203f600a4e8e690c41c1b77611036c8aea96981fef8e85b85c1339cb4a382f2"
$ws.Range("H27").Value = "5002591af40007c5da1e31ae2f2a424f596ecd61ae5ab88d2b8ec064436e550a"
$ws.Range("H28").Value = "c88f5f94121a2839dd49a09cb87c78b2c506d846765cfca6e62656fe73b393fa"
$ws.Range("H29").Value = "9b926bc0adbea331a0a71f7bf18f6c7a62ebde7dd7a52fabe602ad8b00722c56"
$ws.Range("H31").Value = "0543c6777909583f25c96f8942797850af6c7934bdc662b7db680c8d44ed771f"
$ws.Range("H34").Value = "ed47ecd6c9bc82592389aafc247a1e6ae7494f0006ce4205861eed548c773085"

# Update the active selection to match
$ws.Range("C16").Select()
